# LIONBeamLine-Params-Flat.xlsx -- "Trying to implement flexible envelope calculations"
#
# Updates a handful of parameter values on the single worksheet, drops the
# yellow "needs review" highlight from a few rows whose numbers are now
# settled, and leaves the selection on the cell (F4) that was being edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Value updates
# ---------------------------------------------------------------------
$ws.Range("F3").Value = 11
$ws.Range("F11").Value = 0.041180000000000001
$ws.Range("F18").Value = 0.036953
$ws.Range("F22").Value = 1.6

# ---------------------------------------------------------------------
# Formatting clean-up: rows 14, 15, 18, 19 had a yellow "to check" fill
# that is no longer needed now the values are confirmed -- replace it
# with the plain bordered look used elsewhere in the table (matching the
# format already on, e.g., row 11).
# ---------------------------------------------------------------------
$plain = $ws.Range("A11")
$plain.Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Range("G14:H14").PasteSpecial(-4122)
$ws.Range("F14").PasteSpecial(-4122)

$plain.Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)
$ws.Range("G15:H15").PasteSpecial(-4122)
$ws.Range("F15").PasteSpecial(-4122)

$plain.Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)
$ws.Range("G18:H18").PasteSpecial(-4122)
$ws.Range("F18").PasteSpecial(-4122)

$plain.Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$ws.Range("G19:H19").PasteSpecial(-4122)
$ws.Range("F19").PasteSpecial(-4122)

# Row 22 used a slightly different border (no bottom rule); normalise it
# to the same plain bordered style as the rest of the table too.
$plain.Copy()
$ws.Range("A22:E22").PasteSpecial(-4122)
$ws.Range("G22:H22").PasteSpecial(-4122)
$ws.Range("F22").PasteSpecial(-4122)

# F16 picks up the plain border that the rest of its row already has.
$plain.Copy()
$ws.Range("F16").PasteSpecial(-4122)

# F11 and F12 drop their direct formatting entirely (back to the default,
# un-bordered look). ClearFormats only touches formatting, so the values
# set above are left intact.
$ws.Range("F11:F12").ClearFormats()

# ---------------------------------------------------------------------
# Leave the selection where the editor was last working.
# ---------------------------------------------------------------------
$ws.Range("F4").Select()
